# Auto-generated edit script applying the diff's numeric cell updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 16394982
$ws.Cells.Item(112, 10).Value = 1569.9153
$ws.Cells.Item(112, 12).Value = 4709.7459
$ws.Cells.Item(112, 14).Value = -6925.7459

$ws.Cells.Item(137, 8).Value = 3923.0444
$ws.Cells.Item(137, 9).Value = 3720
$ws.Cells.Item(137, 10).Value = 4135.3184
$ws.Cells.Item(137, 11).Value = 11160
$ws.Cells.Item(137, 12).Value = 12405.9552
$ws.Cells.Item(137, 13).Value = -8610
$ws.Cells.Item(137, 14).Value = -17505.9552

$ws.Cells.Item(141, 8).Value = 4221.75
$ws.Cells.Item(141, 9).Value = 4254.3335
$ws.Cells.Item(141, 10).Value = 3863.3333
$ws.Cells.Item(141, 11).Value = 12763.0005
$ws.Cells.Item(141, 12).Value = 11589.9999
$ws.Cells.Item(141, 13).Value = -7583.000499999998
$ws.Cells.Item(141, 14).Value = -21949.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3777.4944
$ws.Cells.Item(32, 9).Value = 3644.9714
$ws.Cells.Item(32, 10).Value = 4219.2383
$ws.Cells.Item(32, 11).Value = 3644.9714
$ws.Cells.Item(32, 12).Value = 4219.2383
$ws.Cells.Item(32, 13).Value = -3357.9714
$ws.Cells.Item(32, 14).Value = -4793.2383

$ws.Cells.Item(45, 8).Value = 836.2917
$ws.Cells.Item(45, 9).Value = 748.2941
$ws.Cells.Item(45, 10).Value = 1050
$ws.Cells.Item(45, 11).Value = 748.2941
$ws.Cells.Item(45, 12).Value = 1050
$ws.Cells.Item(45, 13).Value = -371.2941
$ws.Cells.Item(45, 14).Value = -1804

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1400
$ws.Cells.Item(86, 9).Value = 1400
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 1400
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).Value = -277
$ws.Cells.Item(86, 14).ClearContents()

$ws.Cells.Item(89, 8).Value = 1400
$ws.Cells.Item(89, 9).Value = 1400
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 7000
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).Value = -1384
$ws.Cells.Item(89, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10002984
$ws.Cells.Item(31, 9).Value = 1774.2593
$ws.Cells.Item(31, 10).Value = 21743534
$ws.Cells.Item(31, 11).Value = 1774.2593
$ws.Cells.Item(31, 12).Value = 21743534
$ws.Cells.Item(31, 13).Value = -1479.2593
$ws.Cells.Item(31, 14).Value = -21744124

$ws.Cells.Item(34, 8).Value = 10002984
$ws.Cells.Item(34, 9).Value = 1774.2593
$ws.Cells.Item(34, 10).Value = 21743534
$ws.Cells.Item(34, 11).Value = 1774.2593
$ws.Cells.Item(34, 12).Value = 21743534
$ws.Cells.Item(34, 13).Value = -1572.2593
$ws.Cells.Item(34, 14).Value = -21743938

$ws.Cells.Item(132, 8).Value = 2000.6545
$ws.Cells.Item(132, 9).Value = 1720.7441
$ws.Cells.Item(132, 10).Value = 3003.6667
$ws.Cells.Item(132, 11).Value = 5162.2323
$ws.Cells.Item(132, 12).Value = 9011.000100000001
$ws.Cells.Item(132, 13).Value = -2632.2323
$ws.Cells.Item(132, 14).Value = -14071.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(39, 8).Value = 17571.143
$ws.Cells.Item(39, 10).Value = 17571.143
$ws.Cells.Item(39, 12).Value = 17571.143
$ws.Cells.Item(39, 14).Value = -18635.143

$ws.Cells.Item(41, 8).Value = 12299.125
$ws.Cells.Item(41, 9).Value = 3038.25
$ws.Cells.Item(41, 10).Value = 21560
$ws.Cells.Item(41, 11).Value = 3038.25
$ws.Cells.Item(41, 12).Value = 21560
$ws.Cells.Item(41, 13).Value = -2683.25
$ws.Cells.Item(41, 14).Value = -22270

$ws.Cells.Item(70, 8).Value = 6955.7383
$ws.Cells.Item(70, 9).Value = 5913.7856
$ws.Cells.Item(70, 10).Value = 9039.643
$ws.Cells.Item(70, 11).Value = 5913.7856
$ws.Cells.Item(70, 12).Value = 9039.643
$ws.Cells.Item(70, 13).Value = -5643.7856
$ws.Cells.Item(70, 14).Value = -9579.643

$ws.Cells.Item(73, 8).Value = 6955.7383
$ws.Cells.Item(73, 9).Value = 5913.7856
$ws.Cells.Item(73, 10).Value = 9039.643
$ws.Cells.Item(73, 11).Value = 5913.7856
$ws.Cells.Item(73, 12).Value = 9039.643
$ws.Cells.Item(73, 13).Value = -4977.7856
$ws.Cells.Item(73, 14).Value = -10911.643

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1265.2693
$ws.Cells.Item(16, 9).Value = 1099.7142
$ws.Cells.Item(16, 10).Value = 1960.6
$ws.Cells.Item(16, 11).Value = 1099.7142
$ws.Cells.Item(16, 12).Value = 1960.6
$ws.Cells.Item(16, 13).Value = -929.7141999999999
$ws.Cells.Item(16, 14).Value = -2300.6

$ws.Cells.Item(38, 8).Value = 27996
$ws.Cells.Item(38, 10).Value = 27996
$ws.Cells.Item(38, 12).Value = 27996
$ws.Cells.Item(38, 14).Value = -28816

$ws.Cells.Item(50, 8).Value = 34941.4
$ws.Cells.Item(50, 10).Value = 34941.4
$ws.Cells.Item(50, 12).Value = 34941.4
$ws.Cells.Item(50, 14).Value = -36215.4

$ws.Cells.Item(51, 8).Value = 41996
$ws.Cells.Item(51, 10).Value = 41996
$ws.Cells.Item(51, 12).Value = 41996
$ws.Cells.Item(51, 14).Value = -42952

$ws.Cells.Item(54, 8).Value = 35080.5
$ws.Cells.Item(54, 10).Value = 35080.5
$ws.Cells.Item(54, 12).Value = 35080.5
$ws.Cells.Item(54, 14).Value = -36368.5

$ws.Cells.Item(74, 8).Value = 35773.75
$ws.Cells.Item(74, 10).Value = 39455.715
$ws.Cells.Item(74, 12).Value = 39455.715
$ws.Cells.Item(74, 14).Value = -41451.715

$ws.Cells.Item(77, 8).Value = 35773.75
$ws.Cells.Item(77, 10).Value = 39455.715
$ws.Cells.Item(77, 12).Value = 118367.145
$ws.Cells.Item(77, 14).Value = -128351.145

$ws.Cells.Item(82, 8).Value = 3742.122
$ws.Cells.Item(82, 9).Value = 5901.263
$ws.Cells.Item(82, 11).Value = 5901.263
$ws.Cells.Item(82, 13).Value = -5540.263

$ws.Cells.Item(85, 8).Value = 3742.122
$ws.Cells.Item(85, 9).Value = 5901.263
$ws.Cells.Item(85, 11).Value = 5901.263
$ws.Cells.Item(85, 13).Value = -4653.263

$ws.Cells.Item(94, 8).Value = 34900
$ws.Cells.Item(94, 10).Value = 34900
$ws.Cells.Item(94, 12).Value = 34900
$ws.Cells.Item(94, 14).Value = -36252

$ws.Cells.Item(122, 8).Value = 5281.6875
$ws.Cells.Item(122, 9).Value = 2056.3333
$ws.Cells.Item(122, 10).Value = 9428.571
$ws.Cells.Item(122, 11).Value = 6168.999899999999
$ws.Cells.Item(122, 12).Value = 28285.713
$ws.Cells.Item(122, 13).Value = -3718.999899999999
$ws.Cells.Item(122, 14).Value = -33185.713

$ws.Cells.Item(123, 8).Value = 49710
$ws.Cells.Item(123, 10).Value = 49710
$ws.Cells.Item(123, 12).Value = 49710
$ws.Cells.Item(123, 14).Value = -59510

$ws.Cells.Item(132, 8).Value = 3289.0532
$ws.Cells.Item(132, 9).Value = 1052.4728
$ws.Cells.Item(132, 10).Value = 9439.65
$ws.Cells.Item(132, 11).Value = 3157.4184
$ws.Cells.Item(132, 12).Value = 28318.95
$ws.Cells.Item(132, 13).Value = -627.4184
$ws.Cells.Item(132, 14).Value = -33378.95

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 696.2
$ws.Cells.Item(107, 9).Value = 657
$ws.Cells.Item(107, 10).Value = 741
$ws.Cells.Item(107, 11).Value = 1971
$ws.Cells.Item(107, 12).Value = 2223
$ws.Cells.Item(107, 13).Value = -51
$ws.Cells.Item(107, 14).Value = -6063

$ws.Cells.Item(126, 8).Value = 233328.7
$ws.Cells.Item(126, 9).Value = 1612.619
$ws.Cells.Item(126, 10).Value = 427970.2
$ws.Cells.Item(126, 11).Value = 4837.857
$ws.Cells.Item(126, 12).Value = 1283910.6
$ws.Cells.Item(126, 13).Value = -2367.857
$ws.Cells.Item(126, 14).Value = -1288850.6
